$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update p-values that changed in this revision (block analysis rerun) ---
$ws.Range("E2").Value = 0.8961
$ws.Range("E3").Value = 0.1312
$ws.Range("E4").Value = 0.0002483

# Row 5 ("Chlorophyll Fluorescence (Cold)") no longer differs by North/South site,
# so it switches from the Yes/North-South scientific-notation styling to the
# plain No/"-" styling (same look as rows 3, 6, 7) and gets a new p-value.
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E6").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("C5").Value = "No"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 0.0567

$ws.Range("E6").Value = 0.4441
$ws.Range("E7").Value = 0.6019

# Leave the cursor where the author last left it before saving.
$ws.Range("F8").Select()
